# Add "NA" values under the duplicate_image_filename column (column E),
# for every data row of the first table (rows 2-21), matching the commit
# "add the NA's under duplicate_image_filename".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
